# Fix the row that was left half-updated by the buggy Edit/Delete buttons
# (row 22 should have ended up as id 23 after the edit), and append the
# new todo item that was added afterwards (row 23 / id 24).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 23
$ws.Range("B22").Value = "testeeeeeee"
$ws.Range("C22").Value = "2024-01-26 17:33:50"

$ws.Range("A23").Value = 24
$ws.Range("B23").Value = "Oiiiiiiiii"
$ws.Range("C23").Value = "2024-01-26 17:43:34"
